$wb = $excel.ActiveWorkbook

# --- MicroServiceData sheet: move view / selection away from its old spot ---
$wsUser = $wb.Worksheets.Item("MicroServiceData")
$wsUser.Activate()
$wsUser.Range("AC2").Select()

# --- Pipeline sheet: append a new data row (row 19), copying row 17's layout/format
#     and then overwriting the two cells that actually differ ---
$wsPipeline = $wb.Worksheets.Item("Pipeline")
$wsPipeline.Activate()

$wsPipeline.Range("A17:P17").Copy($wsPipeline.Range("A19:P19"))

$wsPipeline.Range("A19").Value = 1
$wsPipeline.Range("B19").Value = "EditPipeline"
$wsPipeline.Range("C19").Value = "perfeasy-testing"
$wsPipeline.Range("D19").Value = "Version 2"
$wsPipeline.Range("E19").Value = 10
$wsPipeline.Range("F19").Value = "MANUAL"
$wsPipeline.Range("G19").Value = "Build"
$wsPipeline.Range("H19").Value = "devcommunity"
$wsPipeline.Range("I19").Value = "Deploy"
$wsPipeline.Range("J19").Value = "qacommunity"
$wsPipeline.Range("K19").Value = "Generated"
$wsPipeline.Range("L19").Value = "Promote"
$wsPipeline.Range("M19").Value = "Generated"
$wsPipeline.Range("N19").Value = "prodcommunity"
$wsPipeline.Range("O19").Value = "BasicPipeline2be7rzv8g "
$wsPipeline.Range("P19").Value = "CreateBasicPipeLine,CreateJiraPipeLine,runBasicPipeLine"

# View: scroll pipeline sheet and select B20, make Pipeline the active/visible sheet
$wsPipeline.Range("B20").Select()
$wsPipeline.Activate()
